# Generate Report for Handback
#
# - Flips the "Ready for handoff" status string to the handback message
#   (shared string is reused everywhere it already appears: Overview!B2/C2/B3/C3,
#   and the per-language sheets' B2/B3).
# - On each per-language sheet (zh-cn, de-de) adds the "Latest Target File"
#   (col E) and "Latest Handback File" (col F) hyperlinked values for the two
#   localized docs (rows 2 and 3), mirroring the existing handoff file info
#   already present in columns A/C.
# - Stamps "Latest Handback DateTime" (col G) with the real handback
#   timestamp instead of the 0001-01-01 00:00:00 placeholder.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$hyperlinkColor = 15570276   # RGB(100,149,237) == style "FF6495ED" used by the workbook's HyperLink style

function Get-HyperlinkFor($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            return $hl
        }
    }
    return $null
}

function Add-MirrorHyperlink($ws, $srcAddr, $dstAddr) {
    $src = Get-HyperlinkFor $ws $srcAddr
    if ($src -eq $null) {
        return
    }
    $target = $src.Address
    $display = $src.TextToDisplay
    $ws.Hyperlinks.Add($ws.Range($dstAddr), $target, "", "", $display)
    $ws.Range($dstAddr).Font.Underline = $true
    $ws.Range($dstAddr).Font.Color = $hyperlinkColor
}

# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($name -eq "Overview") {
        if ($ws.Range("B2").Value() -eq $oldStatus) { $ws.Range("B2").Value = $newStatus }
        if ($ws.Range("C2").Value() -eq $oldStatus) { $ws.Range("C2").Value = $newStatus }
        if ($ws.Range("B3").Value() -eq $oldStatus) { $ws.Range("B3").Value = $newStatus }
        if ($ws.Range("C3").Value() -eq $oldStatus) { $ws.Range("C3").Value = $newStatus }
    } elseif ($name -eq "zh-cn" -or $name -eq "de-de") {
        if ($ws.Range("B2").Value() -eq $oldStatus) { $ws.Range("B2").Value = $newStatus }
        if ($ws.Range("B3").Value() -eq $oldStatus) { $ws.Range("B3").Value = $newStatus }
    }
}

# 2) Per-language sheets: mirror handoff file info into the Target/Handback
#    File columns, and stamp the real handback datetime.
$ws2 = $wb.Worksheets.Item("zh-cn")
Add-MirrorHyperlink $ws2 '$A$2' 'E2'
Add-MirrorHyperlink $ws2 '$C$2' 'F2'
Add-MirrorHyperlink $ws2 '$A$3' 'E3'
Add-MirrorHyperlink $ws2 '$C$3' 'F3'
$ws2.Range("G2").Value = "2016-03-04 11:19:35"
$ws2.Range("G3").Value = "2016-03-04 11:19:35"

$ws3 = $wb.Worksheets.Item("de-de")
Add-MirrorHyperlink $ws3 '$A$2' 'E2'
Add-MirrorHyperlink $ws3 '$C$2' 'F2'
Add-MirrorHyperlink $ws3 '$A$3' 'E3'
Add-MirrorHyperlink $ws3 '$C$3' 'F3'
$ws3.Range("G2").Value = "2016-03-04 11:19:57"
$ws3.Range("G3").Value = "2016-03-04 11:19:57"

Write-Host "Handback report generated."
